$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D (pushes Sales Area Non Food/Food/Food&Bev right by one).
# The new column inherits the date-formatted style (s="1") from the old column C
# that sits to its left, which is exactly what we want for the new date column.
$ws.Columns.Item(4).Insert()

# Rename header of former "Start Date" column (still column C) -> "Occupancy Start Date"
$ws.Cells.Item(1, 3).Value = "Occupancy Start Date"

# New column D header -> "License Start Date"
$ws.Cells.Item(1, 4).Value = "License Start Date"

# Fill License Start Date values for each data row (as dates, keeping the inherited
# yyyy-mm-dd number format already applied by the column insert)
$ws.Cells.Item(2, 4).Value = "1/1/2008"
$ws.Cells.Item(3, 4).Value = "1/1/2014"
$ws.Cells.Item(4, 4).Value = "1/1/2011"
$ws.Cells.Item(5, 4).Value = "2/15/2020"

# Update the Sales Area Non Food value for OXF-TOPMODEL-001 / OXF-001 (row 5, now col E)
$ws.Cells.Item(5, 5).Value = 201.33

# Column widths: C & D share a wider custom width (not best-fit), matching the new long headers
$ws.Columns.Item(3).ColumnWidth = 18.5
$ws.Columns.Item(4).ColumnWidth = 18.5

# Page setup: portrait, paper size 9 (A4)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Last active selection ends on D5
$ws.Range("D5").Select() | Out-Null
